# Apply updated "想去人数" (interest count) / ticket-price figures across the
# four sheets, matching the regenerated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 8220
$ws.Range("F6").Value  = 103
$ws.Range("F7").Value  = 7167
$ws.Range("F8").Value  = 1137
$ws.Range("F9").Value  = 558
$ws.Range("F10").Value = 496
$ws.Range("F11").Value = 714
$ws.Range("F12").Value = 349
$ws.Range("F18").Value = 11770
$ws.Range("F22").Value = 2330
$ws.Range("F24").Value = 3280
$ws.Range("F27").Value = 2759
$ws.Range("F29").Value = 26
$ws.Range("F31").Value = 3099
$ws.Range("F32").Value = 49
$ws.Range("F33").Value = 2398
$ws.Range("F35").Value = 1640
$ws.Range("F37").Value = 108
$ws.Range("F38").Value = 5873
$ws.Range("F40").Value = 11
$ws.Range("F45").Value = 1082
$ws.Range("F46").Value = 1542
$ws.Range("F47").Value = 9
$ws.Range("F48").Value = 102

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value  = 254
$ws.Range("G16").Value = 180

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 265
$ws.Range("F3").Value = 405

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 265
$ws.Range("F5").Value  = 405
$ws.Range("F8").Value  = 8219
$ws.Range("F10").Value = 103
$ws.Range("F11").Value = 7167
$ws.Range("F12").Value = 7167
$ws.Range("F13").Value = 1137
$ws.Range("F14").Value = 558
$ws.Range("F15").Value = 496
$ws.Range("F16").Value = 714
$ws.Range("F17").Value = 349
$ws.Range("F19").Value = 254
$ws.Range("F22").Value = 11770
$ws.Range("F26").Value = 2330
$ws.Range("F27").Value = 2330
$ws.Range("F28").Value = 3281
$ws.Range("F29").Value = 2759
$ws.Range("F31").Value = 26
$ws.Range("F33").Value = 3103
$ws.Range("F34").Value = 49
$ws.Range("F36").Value = 2398
$ws.Range("F38").Value = 1640
$ws.Range("F39").Value = 108
$ws.Range("F40").Value = 5873
$ws.Range("F48").Value = 1082
$ws.Range("F49").Value = 1542
$ws.Range("F50").Value = 102
